# SystemParameters.xlsx - "Add files via upload" edit
#
# Functional change: the "RZM-350" model line was split into two distinct
# models - a Marine variant (was "RZM-350Marine") and a Land variant
# (was plain "RZM-350"), both renamed with an underscore separator:
#   RZM-350Marine -> RZM-350_Marine   (rows 34-35)
#   RZM-350        -> RZM-350_Land    (rows 30-31)
#
# Column D on every affected row is a shared CONCAT formula
# (=_xlfn.CONCAT(A,"-",B)) so it recalculates automatically once column A
# changes.
#
# NOTE on shared-string ordering: write the Marine rows before the Land
# rows so new shared-string table entries are appended in the same order
# the source workbook shows them (RZM-350_Marine, then RZM-350_Land).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Marine variant (rows 34-35): "RZM-350Marine" -> "RZM-350_Marine"
$ws.Range("A34").Value = "RZM-350_Marine"
$ws.Range("A35").Value = "RZM-350_Marine"

# Land variant (rows 30-31): "RZM-350" -> "RZM-350_Land"
$ws.Range("A30").Value = "RZM-350_Land"
$ws.Range("A31").Value = "RZM-350_Land"

# View state: selection moved to B32, scrolled further down the list.
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
$ws.Range("B32").Select() | Out-Null

# Window geometry / absolute path bookkeeping (best-effort; cosmetic
# metadata written by the authoring machine's Excel instance).
$win.Left = 3000
$win.Top = 1140
$win.Width = 24165
$win.Height = 13590

try {
    $wb.Path = "X:\Python Projects\Optimizer\"
} catch {
}
